$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.163.42'
$ws.Range('E2').Value = '  +0.91%  '
$ws.Range('D3').Value = '2.361.33'
$ws.Range('E3').Value = '  +2.53%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.18'
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.20'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.508'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +0.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.12'
$ws.Range('E10').Value = '  -1.79%  '
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '18.56'
$ws.Range('E12').Value = '  -3.74%  '
$ws.Range('E13').Value = '  +3.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.76'
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('D15').Value = '2.728.41'
$ws.Range('E15').Value = '  +2.69%  '
$ws.Range('D16').Value = '2.357.60'
$ws.Range('E16').Value = '  +1.86%  '
$ws.Range('E17').Value = '  +1.47%  '
$ws.Range('D18').Value = '43.147.36'
$ws.Range('E18').Value = '  +1.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.31'
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('E20').Value = '  +4.19%  '
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.16'
$ws.Range('E22').Value = '  +0.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.73'
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('E24').Value = '  -2.58%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.43'
$ws.Range('E26').Value = '  +0.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.65'
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('E28').Value = '  -0.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.14'
$ws.Range('E29').Value = '  +0.67%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.48'
$ws.Range('E30').Value = '  -2.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.08'
$ws.Range('E32').Value = '  +1.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0727'
$ws.Range('E33').Value = '  +3.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.21'
$ws.Range('E34').Value = '  -2.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.38'
$ws.Range('E35').Value = '  -2.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.83'
$ws.Range('E36').Value = '  +4.63%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.31'
$ws.Range('E37').Value = '  -0.99%  '
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '22.49'
$ws.Range('E39').Value = '  +10.27%  '
$ws.Range('E40').Value = '  +2.19%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '103.44'
$ws.Range('E42').Value = '  -37.35%  '
$ws.Range('D43').Value = '1.940.90'
$ws.Range('E43').Value = '  -1.59%  '
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.13'
$ws.Range('E45').Value = '  +5.52%  '
$ws.Range('E46').Value = '  -9.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.74'
$ws.Range('E47').Value = '  -1.35%  '
$ws.Range('D48').Value = '2.591.82'
$ws.Range('E48').Value = '  +2.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '52.95'
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('E50').Value = '  -0.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.36'
$ws.Range('E51').Value = '  +1.14%  '
